$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 1.036628954199312
$ws.Range("D2").Value = 1.046688860746184
$ws.Range("E2").Value = 1.045267709502108
$ws.Range("F2").Value = 1.055682049324105
$ws.Range("J2").Value = 1.041736231222188
$ws.Range("K2").Value = 1.049453594503636
$ws.Range("L2").Value = 1.048036428021583
$ws.Range("M2").Value = 1.058421836660962
$ws.Range("N2").Value = 1.017722173393022

$ws.Range("C3").Value = 1.037852269132868
$ws.Range("D3").Value = 1.047854453719565
$ws.Range("E3").Value = 1.046375019582421
$ws.Range("F3").Value = 1.056917394539847
$ws.Range("J3").Value = 1.042602211564488
$ws.Range("K3").Value = 1.050429961693895
$ws.Range("L3").Value = 1.048954373455116
$ws.Range("M3").Value = 1.059469596166919
$ws.Range("N3").Value = 1.01802112361675

$ws.Range("C4").Value = 1.038644269418575
$ws.Range("D4").Value = 1.048609449522055
$ws.Range("E4").Value = 1.047092223342602
$ws.Range("F4").Value = 1.057717696286038
$ws.Range("J4").Value = 1.043162496199084
$ws.Range("K4").Value = 1.05106193894649
$ws.Range("L4").Value = 1.049548456834881
$ws.Range("M4").Value = 1.060147951975863
$ws.Range("N4").Value = 1.018214255462612

$ws.Range("C5").Value = 1.038977331835795
$ws.Range("D5").Value = 1.048927037288179
$ws.Range("E5").Value = 1.047393904225826
$ws.Range("F5").Value = 1.058054371844491
$ws.Range("J5").Value = 1.04339802569527
$ws.Range("K5").Value = 1.05132767152284
$ws.Range("L5").Value = 1.04979823691199
$ws.Range("M5").Value = 1.060433225944287
$ws.Range("N5").Value = 1.018295374207536

$ws.Range("C6").Value = 1.039033260716112
$ws.Range("D6").Value = 1.048980372697275
$ws.Range("E6").Value = 1.047444567673858
$ws.Range("F6").Value = 1.05811091459519
$ws.Range("J6").Value = 1.043437571327152
$ws.Range("K6").Value = 1.051372292110849
$ws.Range("L6").Value = 1.04984017771692
$ws.Range("M6").Value = 1.060481130190604
$ws.Range("N6").Value = 1.018308990071483

$ws.Range("C7").Value = 1.038648719398269
$ws.Range("D7").Value = 1.04861369240722
$ws.Range("E7").Value = 1.047096253754515
$ws.Range("F7").Value = 1.057722194061311
$ws.Range("J7").Value = 1.043165643410378
$ws.Range("K7").Value = 1.051065489483198
$ws.Range("L7").Value = 1.049551794299763
$ws.Range("M7").Value = 1.06015176345326
$ws.Range("N7").Value = 1.018215339665466

$ws.Range("C8").Value = 1.037042290401961
$ws.Range("D8").Value = 1.047082617601185
$ws.Range("E8").Value = 1.045641785660282
$ws.Range("F8").Value = 1.056099343807887
$ws.Range("J8").Value = 1.042028906772424
$ws.Range("K8").Value = 1.049783520992138
$ws.Range("L8").Value = 1.048346629399699
$ws.Range("M8").Value = 1.05877585266336
$ws.Range("N8").Value = 1.017823269035325

$ws.Range("C9").Value = 1.03421480480665
$ws.Range("D9").Value = 1.04439059172736
$ws.Range("E9").Value = 1.043084152173089
$ws.Range("F9").Value = 1.053246907630094
$ws.Range("J9").Value = 1.040025316128711
$ws.Range("K9").Value = 1.047526039081714
$ws.Range("L9").Value = 1.046223792402088
$ws.Range("M9").Value = 1.056354231511486
$ws.Range("N9").Value = 1.017130020759942

$ws.Range("C10").Value = 1.032331892453018
$ws.Range("D10").Value = 1.042599826516326
$ws.Range("E10").Value = 1.041382589948672
$ws.Range("F10").Value = 1.051350076752665
$ws.Range("J10").Value = 1.038689189921879
$ws.Range("K10").Value = 1.046022009237241
$ws.Range("L10").Value = 1.044809066756409
$ws.Range("M10").Value = 1.054741713125958
$ws.Range("N10").Value = 1.016666255919269

$ws.Range("C11").Value = 1.031517039518687
$ws.Range("D11").Value = 1.041825319973305
$ws.Range("E11").Value = 1.040646618621679
$ws.Range("F11").Value = 1.050529848675587
$ws.Range("J11").Value = 1.038110525614089
$ws.Range("K11").Value = 1.045370961591956
$ws.Range("L11").Value = 1.044196581433664
$ws.Range("M11").Value = 1.054043910430266
$ws.Range("N11").Value = 1.016465059624069

$ws.Range("C12").Value = 1.03121443390717
$ws.Range("D12").Value = 1.041537768601481
$ws.Range("E12").Value = 1.040373367580796
$ws.Range("F12").Value = 1.050225345062095
$ws.Range("J12").Value = 1.037895565842224
$ws.Range("K12").Value = 1.045129163409233
$ws.Range("L12").Value = 1.04396909107177
$ws.Range("M12").Value = 1.053784778510911
$ws.Range("N12").Value = 1.016390268627863

$ws.Range("C13").Value = 1.031279340838263
$ws.Range("D13").Value = 1.041599443266807
$ws.Range("E13").Value = 1.04043197537818
$ws.Range("F13").Value = 1.050290654650193
$ws.Range("J13").Value = 1.037941676279388
$ws.Range("K13").Value = 1.045181028606875
$ws.Range("L13").Value = 1.044017887912624
$ws.Range("M13").Value = 1.053840360353852
$ws.Range("N13").Value = 1.016406314165964

$ws.Range("C14").Value = 1.031492024689016
$ws.Range("D14").Value = 1.041801548144331
$ws.Range("E14").Value = 1.040624029134005
$ws.Range("F14").Value = 1.050504674941969
$ws.Range("J14").Value = 1.038092757347699
$ws.Range("K14").Value = 1.045350973875878
$ws.Range("L14").Value = 1.044177776730992
$ws.Range("M14").Value = 1.054022489225796
$ws.Range("N14").Value = 1.016458878555808

$ws.Range("C15").Value = 1.031623075049672
$ws.Range("D15").Value = 1.041926089456927
$ws.Range("E15").Value = 1.04074237585382
$ws.Range("F15").Value = 1.050636561808611
$ws.Range("J15").Value = 1.038185840950287
$ws.Range("K15").Value = 1.045455686690419
$ws.Range("L15").Value = 1.044276291325722
$ws.Range("M15").Value = 1.054134713153811
$ws.Range("N15").Value = 1.016491257555439

$ws.Range("C16").Value = 1.03238598107879
$ws.Range("D16").Value = 1.042651246949498
$ws.Range("E16").Value = 1.041431451041526
$ws.Range("F16").Value = 1.051404535895209
$ws.Range("J16").Value = 1.038727591492256
$ws.Range("K16").Value = 1.046065221385821
$ws.Range("L16").Value = 1.044849717407741
$ws.Range("M16").Value = 1.05478803290842
$ws.Range("N16").Value = 1.016679600568474

$ws.Range("C17").Value = 1.032864653661016
$ws.Range("D17").Value = 1.043106360846577
$ws.Range("E17").Value = 1.041863907434799
$ws.Range("F17").Value = 1.051886562472492
$ws.Range("J17").Value = 1.039067386284309
$ws.Range("K17").Value = 1.046447621001324
$ws.Range("L17").Value = 1.045209438493404
$ws.Range("M17").Value = 1.055197956913117
$ws.Range("N17").Value = 1.016797640479036

$ws.Range("C18").Value = 1.033143899860673
$ws.Range("D18").Value = 1.043371908650836
$ws.Range("E18").Value = 1.042116230777755
$ws.Range("F18").Value = 1.052167827833086
$ws.Range("J18").Value = 1.039265572027012
$ws.Range("K18").Value = 1.046670688303306
$ws.Range("L18").Value = 1.045419267461295
$ws.Range("M18").Value = 1.055437100072872
$ws.Range("N18").Value = 1.016866454215009

$ws.Range("C19").Value = 1.033239123224477
$ws.Range("D19").Value = 1.043462468483621
$ws.Range("E19").Value = 1.042202279925504
$ws.Range("F19").Value = 1.052263750316043
$ws.Range("J19").Value = 1.039333146443852
$ws.Range("K19").Value = 1.046746751930742
$ws.Range("L19").Value = 1.045490815402607
$ws.Range("M19").Value = 1.055518648775101
$ws.Range("N19").Value = 1.016889911659239

$ws.Range("C20").Value = 1.03281329202847
$ws.Range("D20").Value = 1.043057522385883
$ws.Range("E20").Value = 1.041817500810577
$ws.Range("F20").Value = 1.051834834484117
$ws.Range("J20").Value = 1.039030930646089
$ws.Range("K20").Value = 1.046406591074809
$ws.Range("L20").Value = 1.045170842826136
$ws.Range("M20").Value = 1.055153971672586
$ws.Range("N20").Value = 1.016784979736949

$ws.Range("C21").Value = 1.031429392762253
$ws.Range("D21").Value = 1.04174202958385
$ws.Range("E21").Value = 1.040567470745678
$ws.Range("F21").Value = 1.05044164674046
$ws.Range("J21").Value = 1.038048268246914
$ws.Range("K21").Value = 1.045300928414386
$ws.Range("L21").Value = 1.044130693082613
$ws.Range("M21").Value = 1.05396885508173
$ws.Range("N21").Value = 1.016443401248345

$ws.Range("C22").Value = 1.03055966723822
$ws.Range("D22").Value = 1.040915706534488
$ws.Range("E22").Value = 1.039782229494367
$ws.Range("F22").Value = 1.049566651776143
$ws.Range("J22").Value = 1.037430324215214
$ws.Range("K22").Value = 1.044605926805364
$ws.Range("L22").Value = 1.043476789869258
$ws.Range("M22").Value = 1.053224089518598
$ws.Range("N22").Value = 1.016228303173801

$ws.Range("C23").Value = 1.031020689237429
$ws.Range("D23").Value = 1.041353682586395
$ws.Range("E23").Value = 1.040198434561876
$ws.Range("F23").Value = 1.050030412899814
$ws.Range("J23").Value = 1.037757918384761
$ws.Range("K23").Value = 1.044974344377577
$ws.Range("L23").Value = 1.04382342909795
$ws.Range("M23").Value = 1.053618869736779
$ws.Range("N23").Value = 1.016342362465158

$ws.Range("C24").Value = 1.032836500001689
$ws.Range("D24").Value = 1.043079590113597
$ws.Range("E24").Value = 1.041838469723028
$ws.Range("F24").Value = 1.051858207803694
$ws.Range("J24").Value = 1.039047403413278
$ws.Range("K24").Value = 1.046425130670559
$ws.Range("L24").Value = 1.045188282516837
$ws.Range("M24").Value = 1.055173846581577
$ws.Range("N24").Value = 1.016790700695544

$ws.Range("C25").Value = 1.034945403633324
$ws.Range("D25").Value = 1.045085849569942
$ws.Range("E25").Value = 1.043744736204954
$ws.Range("F25").Value = 1.053983480794768
$ws.Range("J25").Value = 1.040543358842241
$ws.Range("K25").Value = 1.048109479202824
$ws.Range("L25").Value = 1.046772504994014
$ws.Range("M25").Value = 1.05697994007606
$ws.Range("N25").Value = 1.017309523504505
